$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "31/12/2025"
$ws.Range("C7").Value = 12.0402277339796

$ws.Range("B13").Value = "31/12/2025"
$ws.Range("C13").Value = 13.2028748448604

$ws.Range("B19").Value = "31/12/2025"
$ws.Range("C19").Value = 17.8947933321368
